# Added a kimppa partner to test data
#
# Row 3 holds the existing "KIRKKOÄYRÄÄN ENERGIAKIMPPA" kimppa record; a
# bunch of its cells pick up the sheet's normal formatting (moving off the
# bare default style). A brand-new row 4 is appended underneath it: a
# second partner belonging to that same kimppa, shaped just like row 3 but
# with its own contractor/address details and slightly different counters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Assigning a plain string to a cell lets Excel "smart parse" it (e.g.
# "1.1.2023" silently becomes a date serial). The source data stores these
# verbatim as text, so force text formatting for the write, then drop the
# cell back to General - the normal look for every other text cell on this
# sheet - once the literal string is safely in place.
function Set-TextValue($rng, [string]$val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
}

# --- normalize formatting on row 3 -------------------------------------
# Values are untouched; only the visual formatting of these cells changes
# to match the rest of the row/sheet instead of the bare default style.
$row3FormatCols = @("O3","S3","U3","V3","W3","X3","Y3","AA3","AB3","AK3","AL3","AM3","AN3","AP3")
foreach ($addr in $row3FormatCols) {
    $ws.Range($addr).NumberFormat = "General"
}

# --- append row 4: the new kimppa partner -------------------------------
Set-TextValue $ws.Range("A4") $ws.Range("A3").Text

Set-TextValue $ws.Range("E4") "01-0000123-02"
Set-TextValue $ws.Range("F4") "100456789B"
Set-TextValue $ws.Range("G4") "KUVAKALLIONTIE 1"
Set-TextValue $ws.Range("H4") "15230 LAHTI"

Set-TextValue $ws.Range("K4") "RIKU FORSSTRÖM"
Set-TextValue $ws.Range("L4") "RIKU FORSSTRÖM"
Set-TextValue $ws.Range("M4") "KUVAKALLIONTIE 1"
Set-TextValue $ws.Range("N4") "15230 LAHTI"

Set-TextValue $ws.Range("O4") $ws.Range("O3").Text

Set-TextValue $ws.Range("Q4") $ws.Range("Q3").Text
Set-TextValue $ws.Range("R4") $ws.Range("R3").Text
Set-TextValue $ws.Range("S4") $ws.Range("S3").Text
Set-TextValue $ws.Range("T4") $ws.Range("T3").Text

$ws.Range("U4").Value = 22
$ws.Range("V4").Value = 2
Set-TextValue $ws.Range("W4") $ws.Range("W3").Text
$ws.Range("X4").Value = 7
$ws.Range("Y4").Value = 1

$ws.Range("AA4").Value = 1
$ws.Range("AB4").Value = 53
Set-TextValue $ws.Range("AC4") $ws.Range("AC3").Text
Set-TextValue $ws.Range("AD4") $ws.Range("AD3").Text

Set-TextValue $ws.Range("AK4") $ws.Range("AK3").Text
Set-TextValue $ws.Range("AL4") $ws.Range("AL3").Text
Set-TextValue $ws.Range("AM4") $ws.Range("AM3").Text
Set-TextValue $ws.Range("AN4") $ws.Range("AN3").Text
Set-TextValue $ws.Range("AO4") $ws.Range("AO3").Text
$ws.Range("AP4").Value = 560

# Row 4 should read like row 3 visually too (General format, not the bare
# default style).
$row4FormatCols = @("A4","E4","G4","H4","K4","L4","M4","N4","Q4","R4","S4","T4","W4","AC4","AD4","AK4","AL4","AM4","AN4","AO4")
foreach ($addr in $row4FormatCols) {
    $ws.Range($addr).NumberFormat = "General"
}

# Keep the active selection on the newly-added row, as the source file does.
$ws.Range("A4").Select()
